$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell (F1): "notes" ---
$ws.Cells.Item(1, 6).Value = "notes"

# --- Notes for existing updated rows (introduces the two "updated using..." strings
#     in the order they are first used: row7 text, then row6 text) ---
$ws.Cells.Item(7, 6).Value = "updated using statement of 5th May 23"
$ws.Cells.Item(6, 6).Value = "updated using statement of 5th May 23 - includes EPG"

# --- Rename "elec" -> "elec_imp" throughout column A (rows 2,3,6,7) ---
$ws.Cells.Item(2, 1).Value = "elec_imp"
$ws.Cells.Item(3, 1).Value = "elec_imp"
$ws.Cells.Item(6, 1).Value = "elec_imp"
$ws.Cells.Item(7, 1).Value = "elec_imp"

# --- New "elec_exp" fuel value, used later on row 12 ---
$ws.Cells.Item(12, 1).Value = "elec_exp"

# --- Notes for the brand-new rows 10/13/14 ("Notice 20/6/23") ---
$ws.Cells.Item(10, 6).Value = "Notice 20/6/23"

# --- Note unique to row 11 ---
$ws.Cells.Item(11, 6).Value = "Notice 20/6/23 claimed this was current " + [char]0x00A3 + " but it's not!"

# --- Note unique to row 12 ---
$ws.Cells.Item(12, 6).Value = "statement of 5th May 23"

# === Update existing rows 6-9: add dateEnd, tweak price, copy date style ===

# Row 6: elec_imp / kWh -- dateEnd 45107, price now a formula, note already set above
$ws.Cells.Item(6, 3).Copy()
$ws.Cells.Item(6, 4).PasteSpecial(-4122)
$ws.Cells.Item(6, 4).Value = 45107
$ws.Cells.Item(6, 5).Formula = "=0.4893-0.1657"

# Row 7: elec_imp / sc -- dateEnd 45107, price updated, note already set above
$ws.Cells.Item(7, 3).Copy()
$ws.Cells.Item(7, 4).PasteSpecial(-4122)
$ws.Cells.Item(7, 4).Value = 45107
$ws.Cells.Item(7, 5).Value = 0.4001

# Row 8: gas / kWh -- dateEnd 45107, price now a formula, note
$ws.Cells.Item(8, 3).Copy()
$ws.Cells.Item(8, 4).PasteSpecial(-4122)
$ws.Cells.Item(8, 4).Value = 45107
$ws.Cells.Item(8, 5).Formula = "=0.119-0.0219"
$ws.Cells.Item(8, 6).Value = "updated using statement of 5th May 23 - includes EPG"

# Row 9: gas / sc -- dateEnd 45107, price updated, note
$ws.Cells.Item(9, 3).Copy()
$ws.Cells.Item(9, 4).PasteSpecial(-4122)
$ws.Cells.Item(9, 4).Value = 45107
$ws.Cells.Item(9, 5).Value = 0.2616
$ws.Cells.Item(9, 6).Value = "updated using statement of 5th May 23"

# === New rows 10-14 ===

# Row 10: elec_imp / kWh, dateStart 45108, dateEnd blank (but date-styled), price, note
$ws.Cells.Item(10, 1).Value = "elec_imp"
$ws.Cells.Item(10, 2).Value = "kWh"
$ws.Cells.Item(6, 3).Copy()
$ws.Cells.Item(10, 3).PasteSpecial(-4122)
$ws.Cells.Item(10, 3).Value = 45108
$ws.Cells.Item(6, 4).Copy()
$ws.Cells.Item(10, 4).PasteSpecial(-4122)
$ws.Cells.Item(10, 4).Value = ""
$ws.Cells.Item(10, 5).Value = 0.3072
$ws.Cells.Item(10, 6).Value = "Notice 20/6/23"

# Row 11: elec_imp / sc, dateStart 45108, dateEnd blank, price, note
$ws.Cells.Item(11, 1).Value = "elec_imp"
$ws.Cells.Item(11, 2).Value = "sc"
$ws.Cells.Item(6, 3).Copy()
$ws.Cells.Item(11, 3).PasteSpecial(-4122)
$ws.Cells.Item(11, 3).Value = 45108
$ws.Cells.Item(6, 4).Copy()
$ws.Cells.Item(11, 4).PasteSpecial(-4122)
$ws.Cells.Item(11, 4).Value = ""
$ws.Cells.Item(11, 5).Value = 0.42013
$ws.Cells.Item(11, 6).Value = "Notice 20/6/23 claimed this was current " + [char]0x00A3 + " but it's not!"

# Row 12: elec_exp / kWh, dateStart 44887, dateEnd blank, price, note
$ws.Cells.Item(12, 2).Value = "kWh"
$ws.Cells.Item(6, 3).Copy()
$ws.Cells.Item(12, 3).PasteSpecial(-4122)
$ws.Cells.Item(12, 3).Value = 44887
$ws.Cells.Item(6, 4).Copy()
$ws.Cells.Item(12, 4).PasteSpecial(-4122)
$ws.Cells.Item(12, 4).Value = ""
$ws.Cells.Item(12, 5).Value = 0.1135
$ws.Cells.Item(12, 6).Value = "statement of 5th May 23"

# Row 13: gas / kWh, dateStart 45108, dateEnd blank, price, note
$ws.Cells.Item(13, 1).Value = "gas"
$ws.Cells.Item(13, 2).Value = "kWh"
$ws.Cells.Item(6, 3).Copy()
$ws.Cells.Item(13, 3).PasteSpecial(-4122)
$ws.Cells.Item(13, 3).Value = 45108
$ws.Cells.Item(6, 4).Copy()
$ws.Cells.Item(13, 4).PasteSpecial(-4122)
$ws.Cells.Item(13, 4).Value = ""
$ws.Cells.Item(13, 5).Value = 0.07399
$ws.Cells.Item(13, 6).Value = "Notice 20/6/23"

# Row 14: gas / sc, dateStart 45108, dateEnd blank, price, note
$ws.Cells.Item(14, 1).Value = "gas"
$ws.Cells.Item(14, 2).Value = "sc"
$ws.Cells.Item(6, 3).Copy()
$ws.Cells.Item(14, 3).PasteSpecial(-4122)
$ws.Cells.Item(14, 3).Value = 45108
$ws.Cells.Item(6, 4).Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4122)
$ws.Cells.Item(14, 4).Value = ""
$ws.Cells.Item(14, 5).Value = 0.27468
$ws.Cells.Item(14, 6).Value = "Notice 20/6/23"

# --- Selection / active cell like the saved workbook ---
$ws.Range("H17").Select()
